# Re-style the three data tables (slides 14, 15, 16) from the deck's
# custom "Table_0" style to the built-in "Medium Style 2 - Accent 1"
# table style ({30C081B3-B54C-4C02-A509-87E3BCFC0DB0}).
#
# Table styles can't be assigned through the Table.Style property
# directly (PowerPoint's object model requires ApplyStyle for that), so
# walk every slide/shape, find the tables, and apply the new style to
# any table still using the old GUID.

$oldStyleId = "{50C7EF46-2E12-470A-9EEE-DF136F7E0619}"
$newStyleId = "{30C081B3-B54C-4C02-A509-87E3BCFC0DB0}"

$p = $ppt.ActivePresentation

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $s = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $s.Shapes.Count; $shi++) {
        $sh = $s.Shapes.Item($shi)
        if ($sh.HasTable) {
            $tbl = $sh.Table
            if ($tbl.Style -eq $oldStyleId) {
                $tbl.ApplyStyle($newStyleId)
            }
        }
    }
}
